$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-level updates for the crypto price table: new Price (D), Volume(1h) (E),
# and Hora (G) values. Column G moves 6 -> 7 for every data row (2-51); D/E
# are refreshed only where the source data changed.
$updates = @(
  @{ Row = 2; D = '288.72'; E = '1.33%'; G = '7' },
  @{ Row = 3; D = '29.33'; E = '2.71%'; G = '7' },
  @{ Row = 4; D = '5.076'; E = '2.96%'; G = '7' },
  @{ Row = 5; D = '0.06685'; E = '3.08%'; G = '7' },
  @{ Row = 6; D = '7.338'; E = '1.56%'; G = '7' },
  @{ Row = 7; D = '3.405'; E = '1.11%'; G = '7' },
  @{ Row = 8; D = '1.375'; E = '3.62%'; G = '7' },
  @{ Row = 9; D = '0.9175'; E = '0.52%'; G = '7' },
  @{ Row = 10; E = '2.64%'; G = '7' },
  @{ Row = 11; D = '0.06820'; E = '8.69%'; G = '7' },
  @{ Row = 12; D = '0.07654'; E = '0.54%'; G = '7' },
  @{ Row = 13; D = '0.02935'; E = '-1.77%'; G = '7' },
  @{ Row = 14; D = '0.08992'; E = '0.41%'; G = '7' },
  @{ Row = 15; D = '0.001561'; E = '-2.35%'; G = '7' },
  @{ Row = 16; D = '0.04502'; E = '0.70%'; G = '7' },
  @{ Row = 17; D = '0.0006459'; E = '-1.13%'; G = '7' },
  @{ Row = 18; D = '0.006251'; E = '3.09%'; G = '7' },
  @{ Row = 19; G = '7' },
  @{ Row = 20; D = '2.221'; E = '-0.88%'; G = '7' },
  @{ Row = 21; E = '2.03%'; G = '7' },
  @{ Row = 22; E = '-2.40%'; G = '7' },
  @{ Row = 23; D = '4.065'; E = '2.04%'; G = '7' },
  @{ Row = 24; E = '1.69%'; G = '7' },
  @{ Row = 25; D = '0.001189'; E = '0.02%'; G = '7' },
  @{ Row = 26; D = '0.004119'; E = '-4.70%'; G = '7' },
  @{ Row = 27; E = '1.57%'; G = '7' },
  @{ Row = 28; D = '0.0001616'; E = '-1.22%'; G = '7' },
  @{ Row = 29; G = '7' },
  @{ Row = 30; G = '7' },
  @{ Row = 31; G = '7' },
  @{ Row = 32; G = '7' },
  @{ Row = 33; G = '7' },
  @{ Row = 34; G = '7' },
  @{ Row = 35; G = '7' },
  @{ Row = 36; G = '7' },
  @{ Row = 37; G = '7' },
  @{ Row = 38; G = '7' },
  @{ Row = 39; G = '7' },
  @{ Row = 40; D = '0.04235'; E = '1.65%'; G = '7' },
  @{ Row = 41; D = '0.006732'; E = '0.43%'; G = '7' },
  @{ Row = 42; D = '0.1240'; E = '0.72%'; G = '7' },
  @{ Row = 43; E = '-3.74%'; G = '7' },
  @{ Row = 44; D = '0.01374'; E = '16.37%'; G = '7' },
  @{ Row = 45; D = '0.00005695'; E = '5.72%'; G = '7' },
  @{ Row = 46; D = '1.968'; E = '-3.59%'; G = '7' },
  @{ Row = 47; D = '0.01305'; E = '-29.47%'; G = '7' },
  @{ Row = 48; G = '7' },
  @{ Row = 49; G = '7' },
  @{ Row = 50; G = '7' },
  @{ Row = 51; G = '7' }
)

foreach ($u in $updates) {
    if ($u.ContainsKey('D')) {
        $ws.Cells.Item($u.Row, 4).Value = "'" + $u.D
    }
    if ($u.ContainsKey('E')) {
        $ws.Cells.Item($u.Row, 5).Value = "'" + $u.E
    }
    if ($u.ContainsKey('G')) {
        $ws.Cells.Item($u.Row, 7).Value = "'" + $u.G
    }
}
